$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: copy header-row border/bold/alignment formatting from D1 ---
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 1 (header row) becomes numeric 0,1,0,0 ---
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 0
$ws.Range("E1").Value = 0

# --- Row 2: FE ---
$ws.Range("A2").Value = "FE"
$ws.Range("B2").Value = -0.03
$ws.Range("C2").Value = 0.54
$ws.Range("D2").Value = 0.2
$ws.Range("E2").Value = "'"
$ws.Range("E2").Style = "Normal"

# --- Row 3: FE+Disg ---
$ws.Range("A3").Value = "FE+Disg"
$ws.Range("B3").Value = 0.27
$ws.Range("C3").Value = 0.16
$ws.Range("D3").Value = 0.2
$ws.Range("E3").Value = "'"
$ws.Range("E3").Style = "Normal"

# --- Row 4: FE+Disg+Var ---
$ws.Range("A4").Value = "FE+Disg+Var"
$ws.Range("B4").Value = -0.03
$ws.Range("C4").Value = 0.16
$ws.Range("D4").Value = 0.2
$ws.Range("E4").Value = "'"
$ws.Range("E4").Style = "Normal"
